# Allow payslips to take more than 1 section
# - Add a new Sheet3 (copy of Sheet2), placed after Sheet2
# - Remove the SUM() formulas from the TOTAL rows (16, 33, 50) on every
#   sheet, while keeping the existing cell formatting (values become blank)
# - On Sheet1, the TOTAL cells additionally get switched from the
#   accounting/comma number format to a plain "#,##0.00" number format
# - Sheet3's view is zoomed to 160% with C5 selected

$wb = $excel.ActiveWorkbook

$totalCells = @("C16","F16","I16","C33","F33","I33","C50","F50","I50")

# 1. Duplicate Sheet2 to create Sheet3, inserted right after Sheet2
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item($ws2.Index + 1)
$ws3.Name = "Sheet3"

# 2. Strip the SUM formulas from the TOTAL rows on every sheet (keep formatting)
foreach ($wsName in @("Sheet1","Sheet2","Sheet3")) {
    $ws = $wb.Worksheets.Item($wsName)
    foreach ($addr in $totalCells) {
        $ws.Range($addr).ClearContents()
    }
}

# 3. Sheet1's TOTAL cells get a plain number format instead of accounting format
$ws1 = $wb.Worksheets.Item("Sheet1")
foreach ($addr in $totalCells) {
    $ws1.Range($addr).NumberFormat = "#,##0.00"
}

# 4. Adjust Sheet3's view: zoom 160%, select C5
$ws3.Activate()
$excel.ActiveWindow.Zoom = 160
$ws3.Range("C5").Select()

# restore Sheet1 as the active sheet/tab, matching the original workbook
$ws1.Activate()
$ws1.Range("A3").Select()
